# Brighton_stats.xlsx update
# - Data refresh: match attendance + fixture date/time/day correction
# - Data refresh: player "age" (years-days) columns rolled forward by 3 days
# - Tab bookkeeping: "Matches" tab dropped from the visible tab list, the
#   remaining stat tabs shift up one slot, and a trailing placeholder tab
#   ("Sheet_9") is appended.

$wb = $excel.ActiveWorkbook

function Add-AgeDays {
    param([string]$ageStr, [int]$days)
    $parts = $ageStr -split '-'
    $y = [int]$parts[0]
    $d = [int]$parts[1]
    $d = $d + $days
    if ($d -ge 365) {
        $d = $d - 365
        $y = $y + 1
    }
    return ("{0:d2}-{1:d3}" -f $y, $d)
}

# ---------------------------------------------------------------------------
# 1) "Matches" sheet (physically the 1st tab) - attendance + fixture fix
# ---------------------------------------------------------------------------
$matches = $wb.Worksheets.Item(1)

# Attendance for the Man City match (row 36) was missing, now known
$matches.Range("O36").Value = 52471

# Newcastle Utd fixture (row 43) moved from Sat 2025-05-03 15:00 to Sun 2025-05-04 14:00
$matches.Range("B43").Value = "'2025-05-04"
$matches.Range("C43").Value = "14:00"
$matches.Range("F43").Value = "Sun"

# ---------------------------------------------------------------------------
# 2) Player "age" columns (column D) on every stats tab roll forward 3 days
#    Tabs 2-10 (Standard Stats .. Miscellaneous Stats) all share the same
#    player list in column D; "Standard Stats" and "Playing Time" carry a
#    few extra substitute rows (4-42), the rest stop at row 33.
# ---------------------------------------------------------------------------
$statSheetIndexes = 2..10
foreach ($idx in $statSheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)
    $lastRow = 33
    if ($ws.Range("D42").Value2 -match '^\d{2}-\d{3}$') {
        $lastRow = 42
    }
    for ($r = 4; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 4)
        $val = $cell.Value2
        if ($val -match '^\d{2}-\d{3}$') {
            $cell.Value = Add-AgeDays $val 3
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Tab rename/shift: drop "Matches" from the tab names, shift the other
#    names up one slot, and relabel the trailing tab "Sheet_9".
# ---------------------------------------------------------------------------
$finalNames = @(
    "Standard Stats",
    "Shooting Stats",
    "Passing Stats",
    "Pass Types",
    "Goal & Shot Creation",
    "Defensive Actions",
    "Possession",
    "Playing Time",
    "Miscellaneous Stats",
    "Sheet_9"
)

# Stage through unique temp names first so no two tabs ever collide while
# the rename pass is in flight (e.g. tab 1 "Matches" -> "Standard Stats"
# would otherwise clash with the still-unrenamed tab 2).
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $wb.Worksheets.Item($i).Name = "__tmp_rename_${i}__"
}
for ($i = 1; $i -le $finalNames.Count; $i++) {
    $wb.Worksheets.Item($i).Name = $finalNames[$i - 1]
}
